$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.734.87"
$ws.Range("E2").Value = "  -2.18%  "

$ws.Range("D3").Value = "1.538.87"
$ws.Range("E3").Value = "  -1.77%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").Value = "'289.98"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.3896"
$ws.Range("E7").Value = "  +2.82%  "

$ws.Range("D8").Value = "'0.3166"
$ws.Range("E8").Value = "  -4.34%  "

$ws.Range("D9").Value = "'42.67"
$ws.Range("E9").Value = "  -3.66%  "

$ws.Range("D10").Value = "'0.07167"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").Value = "'1.060"
$ws.Range("E11").Value = "  -7.77%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("D13").Value = "'5.603"
$ws.Range("E13").Value = "  -4.71%  "

$ws.Range("D14").Value = "'18.49"
$ws.Range("E14").Value = "  -8.55%  "

$ws.Range("D15").Value = "'6.591"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("D16").Value = "1.541.96"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").Value = "'0.00001093"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "'0.06530"
$ws.Range("E18").Value = "  -1.54%  "

$ws.Range("D19").Value = "'82.67"
$ws.Range("E19").Value = "  -3.76%  "

$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("D21").Value = "'6.129"
$ws.Range("E21").Value = "  -5.00%  "

$ws.Range("D22").Value = "'15.20"
$ws.Range("E22").Value = "  -6.26%  "

$ws.Range("D23").Value = "'10.88"
$ws.Range("E23").Value = "  -7.69%  "

$ws.Range("D24").Value = "'2.359"
$ws.Range("E24").Value = "  +3.50%  "

$ws.Range("D25").Value = "21.748.89"
$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").Value = "'2.378"
$ws.Range("E26").Value = "  -6.62%  "

$ws.Range("D27").Value = "'144.10"
$ws.Range("E27").Value = "  -4.09%  "

$ws.Range("D28").Value = "'18.30"
$ws.Range("E28").Value = "  -4.80%  "

$ws.Range("D29").Value = "'4.849"
$ws.Range("E29").Value = "  -1.17%  "

$ws.Range("D30").Value = "1.712.36"
$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("D31").Value = "'116.82"
$ws.Range("E31").Value = "  -4.36%  "

$ws.Range("D32").Value = "'0.9728"
$ws.Range("E32").Value = "  -15.45%  "

$ws.Range("D33").Value = "'5.895"
$ws.Range("E33").Value = "  -2.05%  "

$ws.Range("D34").Value = "'0.08191"

$ws.Range("D35").Value = "'8.882"
$ws.Range("E35").Value = "  -5.29%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.06035"
$ws.Range("E36").Value = "  -3.44%  "

$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.488"
$ws.Range("E37").Value = "  -21.65%  "

$ws.Range("D38").Value = "'5.083"

$ws.Range("D39").Value = "'0.02199"
$ws.Range("E39").Value = "  -5.94%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2025"
$ws.Range("E40").Value = "  -6.36%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.178"
$ws.Range("E41").Value = "  -5.44%  "

$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.14%  "

$ws.Range("D43").Value = "'10.54"
$ws.Range("E43").Value = "  -5.41%  "

$ws.Range("D44").Value = "'0.5747"
$ws.Range("E44").Value = "  -5.02%  "

$ws.Range("D45").Value = "'3.734"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").Value = "'12.87"
$ws.Range("E46").Value = "  -6.60%  "

$ws.Range("D47").Value = "'0.5518"
$ws.Range("E47").Value = "  -5.92%  "

$ws.Range("D48").Value = "'116.52"
$ws.Range("E48").Value = "  -4.59%  "

$ws.Range("D49").Value = "'1.871"
$ws.Range("E49").Value = "  -6.43%  "

$ws.Range("E50").Value = "  -4.95%  "

$ws.Range("D51").Value = "'0.06719"
$ws.Range("E51").Value = "  -4.53%  "
